# Updates the crypto tracker's Price (D) and Volume(1h) (E) columns to the
# latest scraped values (GitHub Actions run, Mon Oct 30 09:56:07 UTC 2023).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = '34.582.26'
$ws.Range("E2").Value = '  +0.94%  '
$ws.Range("D3").Value = '1.820.79'
$ws.Range("E3").Value = '  +1.72%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'" + '228.12'
$ws.Range("D6").Value = "'" + '0.561'
$ws.Range("E6").Value = '  +1.85%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").Value = "'" + '34.90'
$ws.Range("E8").Value = '  +8.09%  '
$ws.Range("E9").Value = '  +1.49%  '
$ws.Range("E10").Value = '  +0.60%  '
$ws.Range("E11").Value = '  +0.47%  '
$ws.Range("D12").Value = '2.084.74'
$ws.Range("E12").Value = '  +1.76%  '
$ws.Range("D13").Value = "'" + '11.49'
$ws.Range("E13").Value = '  +3.11%  '
$ws.Range("D14").Value = '1.822.91'
$ws.Range("E14").Value = '  +2.42%  '
$ws.Range("D15").Value = "'" + '0.645'
$ws.Range("E15").Value = '  +3.10%  '
$ws.Range("D16").Value = '34.597.33'
$ws.Range("D17").Value = "'" + '4.35'
$ws.Range("E17").Value = '  +3.61%  '
$ws.Range("D18").Value = "'" + '69.23'
$ws.Range("E18").Value = '  +1.69%  '
$ws.Range("D19").Value = "'" + '247.25'
$ws.Range("E19").Value = '  +0.27%  '
$ws.Range("E20").Value = '  -0.04%  '
$ws.Range("D21").Value = "'" + '11.53'
$ws.Range("E21").Value = '  +5.13%  '
$ws.Range("E23").Value = '  +0.37%  '
$ws.Range("D24").Value = "'" + '171.97'
$ws.Range("E24").Value = '  +6.48%  '
$ws.Range("D25").Value = "'" + '2.09'
$ws.Range("E25").Value = '  +1.43%  '
$ws.Range("D26").Value = "'" + '7.37'
$ws.Range("E26").Value = '  +2.61%  '
$ws.Range("D27").Value = "'" + '16.81'
$ws.Range("E27").Value = '  +2.91%  '
$ws.Range("D28").Value = "'" + '0.117'
$ws.Range("E28").Value = '  +1.65%  '
$ws.Range("E29").Value = '  -0.39%  '
$ws.Range("D30").Value = "'" + '4.02'
$ws.Range("E30").Value = '  +6.97%  '
$ws.Range("D31").Value = "'" + '0.0531'
$ws.Range("E31").Value = '  +1.86%  '
$ws.Range("D32").Value = "'" + '3.85'
$ws.Range("E32").Value = '  +2.70%  '
$ws.Range("E33").Value = '  +1.14%  '
$ws.Range("E34").Value = '  +2.40%  '
$ws.Range("E35").Value = '  +1.46%  '
$ws.Range("D36").Value = '1.416.46'
$ws.Range("E36").Value = '  -1.70%  '
$ws.Range("D37").Value = "'" + '0.679'
$ws.Range("E37").Value = '  +2.00%  '
$ws.Range("E38").Value = '  +1.36%  '
$ws.Range("D39").Value = "'" + '86.25'
$ws.Range("E39").Value = '  +5.19%  '
$ws.Range("E40").Value = '  +0.36%  '
$ws.Range("D41").Value = "'" + '2.86'
$ws.Range("E41").Value = '  +4.61%  '
$ws.Range("D42").Value = "'" + '0.955'
$ws.Range("E42").Value = '  +3.63%  '
$ws.Range("E43").Value = '  +0.98%  '
$ws.Range("D44").Value = "'" + '14.00'
$ws.Range("E44").Value = '  -0.59%  '
$ws.Range("D45").Value = "'" + '0.0525'
$ws.Range("E45").Value = '  +1.09%  '
$ws.Range("E46").Value = '  +3.13%  '
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("D48").Value = '1.985.26'
$ws.Range("E48").Value = '  +2.05%  '
$ws.Range("D49").Value = "'" + '105.97'
$ws.Range("E49").Value = '  +0.42%  '
$ws.Range("E50").Value = '  +1.68%  '
$ws.Range("D51").Value = "'" + '1.00'
$ws.Range("E51").Value = '  -0.11%  '

Write-Host "Updated 79 cells."
